{"js": "// Task 2 report edit:\n//  1) \"Constructor method for\" -> \"Constructor method \" + \"Player()\" (Code\n//     character style) + \" for\"\n//  2) \"...hard coded within the method. \" -> \"...hard coded within the\n//     method as per the instruction of the related comment.\"\n\nconst body = context.document.body;\n\n// --- Part 1 ---------------------------------------------------------\n// Replace the original text with a version that keeps a unique\n// placeholder where the code-styled \"Player()\" needs to go, so we can\n// re-search for just that placeholder and style it without the newly\n// inserted text inheriting formatting from its styled neighbour.\nconst constructorMatches = body.search(\"Constructor method for\", { matchCase: true });\nconstructorMatches.load(\"text\");\nawait context.sync();\n\nif (constructorMatches.items.length === 0) {\n  throw new Error('Could not find \"Constructor method for\" in the document.');\n}\n\nconstructorMatches.items[0].insertText(\n  \"Constructor method \\u0001PLAYERCALL\\u0001 for\",\n  \"Replace\"\n);\nawait context.sync();\n\nconst placeholderMatches = body.search(\"\\u0001PLAYERCALL\\u0001\", { matchCase: true });\nplaceholderMatches.load(\"text\");\nawait context.sync();\n\nif (placeholderMatches.items.length === 0) {\n  throw new Error(\"Could not find the Player() placeholder after inserting it.\");\n}\n\nconst playerRange = placeholderMatches.items[0].insertText(\"Player()\", \"Replace\");\nplayerRange.style = \"CodeChar\";\nawait context.sync();\n\n// --- Part 2 ---------------------------------------------------------\nconst paramMatches = body.search(\n  \"The parameters for the constructor are hard coded within the method. \",\n  { matchCase: true }\n);\nparamMatches.load(\"text\");\nawait context.sync();\n\nif (paramMatches.items.length === 0) {\n  throw new Error(\"Could not find the constructor-parameters sentence.\");\n}\n\nparamMatches.items[0].insertText(\n  \"The parameters for the constructor are hard coded within the method\" +\n    \" as per the instruction of the related comment.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Task 2 report edit:\n#  1) \"Constructor method for\" -> \"Constructor method \" + \"Player()\" (Code\n#     character style) + \" for\"\n#  2) \"...hard coded within the method. \" -> \"...hard coded within the\n#     method as per the instruction of the related comment.\"\n\n$d = $word.ActiveDocument\n\n# --- Part 1 ---------------------------------------------------------\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.ClearFormatting()\n$find1.Text = \"Constructor method for\"\n$found1 = $find1.Execute()\nif (-not $found1) {\n    throw 'Could not find \"Constructor method for\" in the document.'\n}\n$range1.Text = \"Constructor method Player() for\"\n\n# Re-find just the freshly inserted \"Player()\" text and apply the Code\n# character style to it (re-searching keeps the surrounding plain-text\n# runs from inheriting the style).\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Player()\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n    throw \"Could not find the Player() placeholder after inserting it.\"\n}\n$range2.Style = \"CodeChar\"\n\n# --- Part 2 ---------------------------------------------------------\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.ClearFormatting()\n$find3.Text = \"The parameters for the constructor are hard coded within the method. \"\n$found3 = $find3.Execute()\nif (-not $found3) {\n    throw \"Could not find the constructor-parameters sentence.\"\n}\n$range3.Text = \"The parameters for the constructor are hard coded within the method as per the instruction of the related comment.\"\n"}
